# Auto-generated edit script: update scraped market-price columns (H-N)
# per sheet (Kraken_Profits workbook: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 773.5909
$ws.Range("I6").Value = 450.58334
$ws.Range("J6").Value = 1161.2
$ws.Range("K6").Value = 1351.75002
$ws.Range("L6").Value = 3483.6
$ws.Range("M6").Value = -1239.75002
$ws.Range("N6").Value = -3707.6
$ws.Range("H7").Value = 1556.25
$ws.Range("I7").Value = 687.5
$ws.Range("J7").Value = 2425
$ws.Range("K7").Value = 687.5
$ws.Range("L7").Value = 2425
$ws.Range("M7").Value = -575.5
$ws.Range("N7").Value = -2649
$ws.Range("H9").Value = 137.5
$ws.Range("I9").Value = 175
$ws.Range("K9").Value = 175
$ws.Range("M9").Value = -6
$ws.Range("H10").Value = 1601.25
$ws.Range("I10").Value = 2250
$ws.Range("J10").Value = 952.5
$ws.Range("K10").Value = 2250
$ws.Range("L10").Value = 952.5
$ws.Range("M10").Value = -1957
$ws.Range("N10").Value = -1538.5
$ws.Range("H14").Value = 1556.25
$ws.Range("I14").Value = 687.5
$ws.Range("J14").Value = 2425
$ws.Range("K14").Value = 687.5
$ws.Range("L14").Value = 2425
$ws.Range("M14").Value = -496.5
$ws.Range("N14").Value = -2807
$ws.Range("H18").Value = 19333.334
$ws.Range("J18").Value = 12200
$ws.Range("L18").Value = 12200
$ws.Range("N18").Value = -12768
$ws.Range("H40").Value = 7567.8623
$ws.Range("I40").Value = 3899
$ws.Range("J40").Value = 8332.208
$ws.Range("K40").Value = 3899
$ws.Range("L40").Value = 8332.208
$ws.Range("M40").Value = -3724
$ws.Range("N40").Value = -8682.208
$ws.Range("H53").Value = 159.85715
$ws.Range("I53").Value = 141.22223
$ws.Range("J53").Value = 193.4
$ws.Range("K53").Value = 141.22223
$ws.Range("L53").Value = 193.4
$ws.Range("M53").Value = 495.77777
$ws.Range("N53").Value = -1467.4
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H64").Value = 33336000
$ws.Range("I64").Value = 33336000
$ws.Range("K64").Value = 33336000
$ws.Range("M64").Value = -33335752
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H67").Value = 33336000
$ws.Range("I67").Value = 33336000
$ws.Range("K67").Value = 33336000
$ws.Range("M67").Value = -33335142
$ws.Range("H69").Value = 6000
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H70").Value = 23899.4
$ws.Range("J70").Value = 36499
$ws.Range("L70").Value = 109497
$ws.Range("N70").Value = -110037
$ws.Range("H72").Value = 6000
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H73").Value = 23899.4
$ws.Range("J73").Value = 36499
$ws.Range("L73").Value = 109497
$ws.Range("N73").Value = -111369

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3224.75
$ws.Range("I61").Value = 1966.3334
$ws.Range("K61").Value = 1966.3334
$ws.Range("M61").Value = -1754.3334
$ws.Range("H86").Value = 5000
$ws.Range("J86").Value = 5000
$ws.Range("L86").Value = 5000
$ws.Range("N86").Value = -7372
$ws.Range("H89").Value = 5000
$ws.Range("J89").Value = 5000
$ws.Range("L89").Value = 15000
$ws.Range("N89").Value = -26856
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").ClearContents()
$ws.Range("N96").Value = 0
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("N113").Value = 0
$ws.Range("H132").Value = 1415
$ws.Range("J132").Value = 1350
$ws.Range("L132").Value = 4050
$ws.Range("N132").Value = -9110
$ws.Range("H136").Value = 3224.75
$ws.Range("I136").Value = 1966.3334
$ws.Range("K136").Value = 5899.0002
$ws.Range("M136").Value = -3349.0002

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1571
$ws.Range("I20").Value = 1487.8
$ws.Range("K20").Value = 1487.8
$ws.Range("M20").Value = -1240.8
$ws.Range("H134").Value = 3292.923
$ws.Range("I134").Value = 3292.923
$ws.Range("K134").Value = 9878.769
$ws.Range("M134").Value = -7343.769

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 230.8125
$ws.Range("I7").Value = 278.45456
$ws.Range("J7").Value = 126
$ws.Range("K7").Value = 278.45456
$ws.Range("L7").Value = 126
$ws.Range("M7").Value = -165.45456
$ws.Range("N7").Value = -352
$ws.Range("H22").Value = 516
$ws.Range("I22").Value = 495
$ws.Range("K22").Value = 495
$ws.Range("M22").Value = -145
$ws.Range("H31").Value = 2225.2354
$ws.Range("I31").Value = 1232.125
$ws.Range("J31").Value = 3108
$ws.Range("K31").Value = 1232.125
$ws.Range("L31").Value = 3108
$ws.Range("M31").Value = -937.125
$ws.Range("N31").Value = -3698
$ws.Range("H34").Value = 2225.2354
$ws.Range("I34").Value = 1232.125
$ws.Range("J34").Value = 3108
$ws.Range("K34").Value = 1232.125
$ws.Range("L34").Value = 3108
$ws.Range("M34").Value = -1030.125
$ws.Range("N34").Value = -3512
$ws.Range("H58").Value = 999.6667
$ws.Range("I58").Value = 999.5
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 999.5
$ws.Range("L58").Value = 1000
$ws.Range("M58").Value = -796.5
$ws.Range("N58").Value = -1406
$ws.Range("H100").Value = 99995
$ws.Range("J100").Value = 99995
$ws.Range("L100").Value = 99995
$ws.Range("N100").Value = -102159
$ws.Range("H132").Value = 5000
$ws.Range("J132").Value = 5000
$ws.Range("L132").Value = 15000
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 999.6667
$ws.Range("I136").Value = 999.5
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 2998.5
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -448.5
$ws.Range("N136").Value = -8100

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 52
$ws.Range("I2").Value = 37.4
$ws.Range("J2").Value = 68.22222
$ws.Range("K2").Value = 224.4
$ws.Range("L2").Value = 409.33332
$ws.Range("M2").Value = -111.4
$ws.Range("N2").Value = -635.33332
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H26").Value = 145
$ws.Range("I26").Value = 145
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 435
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -147

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 91.333336
$ws.Range("I2").Value = 123
$ws.Range("J2").Value = 59.666668
$ws.Range("K2").Value = 123
$ws.Range("L2").Value = 59.666668
$ws.Range("M2").Value = -10
$ws.Range("N2").Value = -285.666668
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("M4").Value = 500
$ws.Range("N4").Value = -724
$ws.Range("H123").Value = 69326
$ws.Range("J123").Value = 69326
$ws.Range("L123").Value = 69326
$ws.Range("N123").Value = -74226
$ws.Range("H132").Value = 7716.7144
$ws.Range("I132").Value = 6255.75
$ws.Range("K132").Value = 18767.25
$ws.Range("M132").Value = -16237.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").ClearContents()
$ws.Range("N27").Value = 0
$ws.Range("H45").Value = 44888
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 44888
$ws.Range("K45").Value = 0
$ws.Range("L45").ClearContents()
$ws.Range("M45").Value = 44888
$ws.Range("N45").Value = -45870
$ws.Range("H62").Value = 3333.6667
$ws.Range("J62").Value = 4000
$ws.Range("L62").Value = 4000
$ws.Range("N62").Value = -5248
$ws.Range("H64").Value = 80000
$ws.Range("I64").Value = 80000
$ws.Range("K64").Value = 80000
$ws.Range("M64").Value = -79752
$ws.Range("H65").Value = 3333.6667
$ws.Range("J65").Value = 4000
$ws.Range("L65").Value = 20000
$ws.Range("N65").Value = -26240
$ws.Range("H67").Value = 80000
$ws.Range("I67").Value = 80000
$ws.Range("K67").Value = 80000
$ws.Range("M67").Value = -79142
$ws.Range("H81").Value = 1983.5
$ws.Range("I81").Value = 1983.5
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 3967
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -2906
$ws.Range("H84").Value = 1983.5
$ws.Range("I84").Value = 1983.5
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 19835
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -14531

